$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '301.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-3.05%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.31'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.32%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.066'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.20%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07916'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-3.03%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.887'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-9.72%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.769'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-2.15%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9288'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.26%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1371'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '30.84%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1897'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.03%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09207'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.64%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03434'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-5.84%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09893'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.28%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001428'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.12%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005890'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.22%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.531'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.59%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.049'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.91%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.926'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.35%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3409'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.46%'
$ws.Range("B20").Value = 'ProBitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1301'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.00%'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.053'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.83%'
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2398'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '8.26%'
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04497'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.24%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001214'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.02%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004761'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.47%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001231'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-1.52%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003004'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-32.54%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01850'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-6.20%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04755'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.00%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007340'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-3.07%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009640'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '7.46%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1325'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.16%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002113'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-2.70%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01102'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-5.49%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006254'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-5.45%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.01%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.68'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '8.18%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '10.50%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.01%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002002'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.01%'
